$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "bleu" -> "noir" in statut_label column (B2:B4, B6:B7)
$ws.Range("B2").Value = "noir"
$ws.Range("B3").Value = "noir"
$ws.Range("B4").Value = "noir"
$ws.Range("B6").Value = "noir"
$ws.Range("B7").Value = "noir"

# Replace "pas de résultat ni de publication" -> "pas de résultat postés ni publiés" in statut_name column
$ws.Range("C2").Value = "pas de résultat postés ni publiés"
$ws.Range("C3").Value = "pas de résultat postés ni publiés"
$ws.Range("C4").Value = "pas de résultat postés ni publiés"
$ws.Range("C6").Value = "pas de résultat postés ni publiés"
$ws.Range("C7").Value = "pas de résultat postés ni publiés"

# Replace "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
$ws.Range("C5").Value = "résultat postés ou publiés dans les 36 mois"
